$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.960869565217391,0.0268115942028985,0.794927536231884,0.426811594202899,0.938405797101449,0.00217391304347826,0.973913043478261,0.955797101449275,0.0108695652173913,0.0108695652173913,0.0920289855072464,0.00507246376811594,0.00434782608695652,0.951449275362319,0.201449275362319,0.426811594202899,0.102898550724638,0.00434782608695652,0.022463768115942,0.00507246376811594,0.952173913043478,0.564492753623188,0.0528985507246377)
    3 = @(0.0036231884057971,0.0181159420289855,0.0231884057971014,0.0123188405797101,0.0072463768115942,0.013768115942029,0.00942028985507246,0.0072463768115942,0.0036231884057971,0.0108695652173913,0.0289855072463768,0.144927536231884,0.0036231884057971,0,0.146376811594203,0.0369565217391304,0.0297101449275362,0.00797101449275362,0.00289855072463768,0.00797101449275362,0.0173913043478261,0.0123188405797101,0.0159420289855072)
    4 = @(0.0253623188405797,0.938405797101449,0.173188405797101,0.543478260869565,0.0376811594202899,0.981884057971015,0.00652173913043478,0.0268115942028985,0.96304347826087,0.0557971014492754,0.872463768115942,0.00869565217391304,0.0108695652173913,0.041304347826087,0.648550724637681,0.00217391304347826,0.865217391304348,0.00144927536231884,0.969565217391304,0.986231884057971,0.013768115942029,0.406521739130435,0.913768115942029)
    5 = @(0.0101449275362319,0.0159420289855072,0.00869565217391304,0.0159420289855072,0.0166666666666667,0.00217391304347826,0.0101449275362319,0.0101449275362319,0.022463768115942,0.922463768115942,0.00579710144927536,0.839855072463768,0.981159420289855,0.00652173913043478,0.00289855072463768,0.534057971014493,0.00217391304347826,0.986231884057971,0.00434782608695652,0.00072463768115942,0.0166666666666667,0.0166666666666667,0.0166666666666667)
}

# Columns B..X correspond to indices 2..24
foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $colNum = $i + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $values[$i]
    }
}
